$d = $word.ActiveDocument

function New-WordOpenXmlPackage($innerParagraphXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$apos = [char]0x2019

# --- Edit 1: Title paragraph "KIRCHHOFF'S LAW" -> centered, bold, size 80, single underline ---
$titlePara = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:spacing w:lineRule="auto" w:line="360"/><w:jc w:val="center"/><w:rPr><w:b/><w:b/><w:bCs/><w:sz w:val="80"/><w:szCs w:val="80"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b/><w:bCs/><w:sz w:val="80"/><w:szCs w:val="80"/><w:u w:val="single"/></w:rPr><w:t>KIRCHHOFF' + $apos + 'S LAW</w:t></w:r></w:p>'

$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML((New-WordOpenXmlPackage $titlePara))

# --- Edit 2: merge "vii" + "i" + ". Find the sum..." runs into a single "viii. Find the sum..." run ---
$para17 = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:spacing w:lineRule="auto" w:line="360"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>viii. Find the sum of voltages in each possible loop. After find the sum of voltages</w:t></w:r></w:p>'

$p17 = $d.Paragraphs.Item(17)
$p17.Range.InsertXML((New-WordOpenXmlPackage $para17))

# --- Edit 3: merge "ix" + ". Next you'll have..." runs into a single "ix. Next you'll have..." run ---
$para18 = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:spacing w:lineRule="auto" w:line="360"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Liberation Sans" w:hAnsi="Liberation Sans"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/></w:rPr><w:t>ix. Next you' + $apos + 'll have the sum of voltages. From the junction, you can form an equation for the current coming into and leaving the junction</w:t></w:r></w:p>'

$p18 = $d.Paragraphs.Item(18)
$p18.Range.InsertXML((New-WordOpenXmlPackage $para18))
